# Updated cryptos list on Tue Nov 14 09:21:29 UTC 2023 with GitHub Actions
#
# This refreshes the "cryptos" worksheet with the latest scraped prices /
# 1h volume deltas, and re-orders a handful of coin rows (27/28, 36/37/38,
# 43/44, 51) to reflect their new rank position, matching the upstream
# data refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving its "plain text" nature.
# Several Price column values look like plain decimal numbers (e.g. "247.29")
# which Excel would otherwise auto-convert to a numeric type; prefixing with
# a single quote forces text entry (as in the source sheet), and resetting
# the style back to Normal avoids leaving a stray quote-prefixed number
# format applied to the cell.
function Set-TextValue($cell, $text) {
    $ws.Range($cell).Value = "'" + $text
    $ws.Range($cell).Style = "Normal"
}


$ws.Range("D2").Value = '36.663.55'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '2.057.83'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue "D5" '247.29'
$ws.Range("E5").Value = '  +0.34%  '
Set-TextValue "D6" '0.666'
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E7").Value = '  +0.03%  '
Set-TextValue "D8" '54.98'
$ws.Range("E8").Value = '  -6.78%  '
Set-TextValue "D9" '60.62'
$ws.Range("E9").Value = '  +1.80%  '
Set-TextValue "D10" '0.368'
$ws.Range("E10").Value = '  -2.58%  '
Set-TextValue "D11" '0.0756'
$ws.Range("E11").Value = '  -2.40%  '
$ws.Range("E12").Value = '  -3.03%  '
Set-TextValue "D13" '0.978'
$ws.Range("E13").Value = '  +10.50%  '
Set-TextValue "D14" '14.84'
$ws.Range("E14").Value = '  -3.98%  '
$ws.Range("D15").Value = '2.361.97'
$ws.Range("E15").Value = '  +0.22%  '
Set-TextValue "D16" '5.48'
$ws.Range("E16").Value = '  -4.07%  '
$ws.Range("D17").Value = '2.049.04'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '36.610.17'
$ws.Range("E18").Value = '  -1.05%  '
Set-TextValue "D19" '17.41'
$ws.Range("E19").Value = '  -4.14%  '
Set-TextValue "D20" '72.30'
$ws.Range("E20").Value = '  -2.33%  '
$ws.Range("D21").Value = '0.0₃0863'
$ws.Range("E21").Value = '  -3.22%  '
Set-TextValue "D22" '238.30'
$ws.Range("E22").Value = '  +0.06%  '
Set-TextValue "D23" '5.27'
$ws.Range("E23").Value = '  -3.05%  '
Set-TextValue "D24" '0.999'
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  -2.56%  '
Set-TextValue "D26" '2.33'
$ws.Range("E26").Value = '  +7.82%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D27" '166.41'
$ws.Range("E27").Value = '  -1.80%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D28" '9.28'
$ws.Range("E28").Value = '  -7.30%  '
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("E31").Value = '  +9.06%  '
Set-TextValue "D32" '5.09'
$ws.Range("E32").Value = '  -5.99%  '
Set-TextValue "D33" '4.52'
$ws.Range("E33").Value = '  -3.58%  '
Set-TextValue "D34" '0.0597'
$ws.Range("E34").Value = '  -3.33%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D36" '0.0862'
$ws.Range("E36").Value = '  +2.54%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D37" '2.28'
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D38" '1.84'
$ws.Range("E38").Value = '  +0.36%  '
Set-TextValue "D39" '5.08'
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("E40").Value = '  -5.57%  '
$ws.Range("E41").Value = '  -5.44%  '
Set-TextValue "D42" '0.0216'
$ws.Range("E42").Value = '  -3.58%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D43" '95.19'
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D44" '1.11'
$ws.Range("E44").Value = '  -4.70%  '
Set-TextValue "D45" '0.0921'
$ws.Range("E45").Value = '  -4.27%  '
$ws.Range("D46").Value = '1.416.21'
$ws.Range("E46").Value = '  +8.54%  '
Set-TextValue "D47" '7.61'
$ws.Range("E47").Value = '  +11.27%  '
Set-TextValue "D48" '16.00'
$ws.Range("E48").Value = '  -5.93%  '
$ws.Range("E49").Value = '  +2.38%  '
Set-TextValue "D50" '2.28'
$ws.Range("E50").Value = '  -3.64%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue "D51" '46.01'
$ws.Range("E51").Value = '  +3.42%  '
